# Apply crypto price/volume updates per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.744.64"
$ws.Cells.Item(2, 5).Value = "  +0.38%  "

$ws.Cells.Item(3, 4).Value = "1.603.27"
$ws.Cells.Item(3, 5).Value = "  +0.42%  "

$ws.Cells.Item(4, 5).Value = "  +0.21%  "

$ws.Cells.Item(5, 4).Value = "'211.91"
$ws.Cells.Item(5, 5).Value = "  +0.28%  "

$ws.Cells.Item(6, 5).Value = "  +0.29%  "

$ws.Cells.Item(7, 5).Value = "  +0.24%  "

$ws.Cells.Item(8, 4).Value = "'0.0619"
$ws.Cells.Item(8, 5).Value = "  +0.27%  "

$ws.Cells.Item(9, 4).Value = "'0.247"
$ws.Cells.Item(9, 5).Value = "  +0.45%  "

$ws.Cells.Item(10, 4).Value = "'19.70"
$ws.Cells.Item(10, 5).Value = "  +1.19%  "

$ws.Cells.Item(11, 4).Value = "'0.0846"
$ws.Cells.Item(11, 5).Value = "  +0.67%  "

$ws.Cells.Item(12, 4).Value = "1.829.19"
$ws.Cells.Item(12, 5).Value = "  +0.41%  "

$ws.Cells.Item(13, 4).Value = "1.606.84"
$ws.Cells.Item(13, 5).Value = "  +1.19%  "

$ws.Cells.Item(14, 4).Value = "'4.07"
$ws.Cells.Item(14, 5).Value = "  +1.22%  "

$ws.Cells.Item(15, 4).Value = "'0.525"
$ws.Cells.Item(15, 5).Value = "  +0.51%  "

$ws.Cells.Item(16, 4).Value = "'65.05"
$ws.Cells.Item(16, 5).Value = "  +0.06%  "

$ws.Cells.Item(17, 4).Value = "0.0₃0744"
$ws.Cells.Item(17, 5).Value = "  +0.78%  "

$ws.Cells.Item(18, 4).Value = "'209.83"
$ws.Cells.Item(18, 5).Value = "  +0.38%  "

$ws.Cells.Item(19, 5).Value = "  +0.22%  "

$ws.Cells.Item(20, 4).Value = "'7.14"
$ws.Cells.Item(20, 5).Value = "  +1.43%  "

$ws.Cells.Item(21, 4).Value = "'4.30"
$ws.Cells.Item(21, 5).Value = "  +0.58%  "

$ws.Cells.Item(22, 5).Value = "  -4.50%  "

$ws.Cells.Item(23, 4).Value = "'9.07"

$ws.Cells.Item(24, 4).Value = "'143.80"
$ws.Cells.Item(24, 5).Value = "  -0.03%  "

$ws.Cells.Item(25, 4).Value = "'1.00"
$ws.Cells.Item(25, 5).Value = "  +0.16%  "

$ws.Cells.Item(26, 5).Value = "  -0.28%  "

$ws.Cells.Item(27, 5).Value = "  -0.01%  "

$ws.Cells.Item(28, 4).Value = "'15.37"
$ws.Cells.Item(28, 5).Value = "  +0.47%  "

$ws.Cells.Item(29, 4).Value = "'0.0509"
$ws.Cells.Item(29, 5).Value = "  -0.85%  "

$ws.Cells.Item(30, 4).Value = "'1.15"
$ws.Cells.Item(30, 5).Value = "  +0.07%  "

$ws.Cells.Item(31, 5).Value = "  +1.23%  "

$ws.Cells.Item(32, 5).Value = "  +0.64%  "

$ws.Cells.Item(33, 4).Value = "1.290.96"
$ws.Cells.Item(33, 5).Value = "  +0.14%  "

$ws.Cells.Item(34, 5).Value = "  +1.25%  "

$ws.Cells.Item(35, 5).Value = "  +21.30%  "

$ws.Cells.Item(36, 5).Value = "  +0.33%  "

$ws.Cells.Item(37, 4).Value = "'0.594"
$ws.Cells.Item(37, 5).Value = "  -3.82%  "

$ws.Cells.Item(38, 5).Value = "  -0.35%  "

$ws.Cells.Item(39, 4).Value = "'0.830"
$ws.Cells.Item(39, 5).Value = "  -0.02%  "

$ws.Cells.Item(40, 5).Value = "  -0.33%  "

$ws.Cells.Item(41, 5).Value = "  -0.30%  "

$ws.Cells.Item(42, 5).Value = "  -0.23%  "

$ws.Cells.Item(43, 4).Value = "'63.12"
$ws.Cells.Item(43, 5).Value = "  -0.14%  "

$ws.Cells.Item(44, 4).Value = "1.740.94"
$ws.Cells.Item(44, 5).Value = "  +0.48%  "

$ws.Cells.Item(45, 4).Value = "'90.57"
$ws.Cells.Item(45, 5).Value = "  -0.58%  "

$ws.Cells.Item(46, 2).Value = "RenderToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(46, 4).Value = "'1.56"
$ws.Cells.Item(46, 5).Value = "  -0.37%  "

$ws.Cells.Item(47, 2).Value = "Aptos"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(47, 4).Value = "'6.22"
$ws.Cells.Item(47, 5).Value = "  +21.50%  "

$ws.Cells.Item(48, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(48, 4).Value = "0.0₆0103"
$ws.Cells.Item(48, 5).Value = "  -3.47%  "

$ws.Cells.Item(49, 2).Value = "Algorand"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(49, 4).Value = "'0.102"
$ws.Cells.Item(49, 5).Value = "  +1.25%  "

$ws.Cells.Item(50, 4).Value = "'0.0514"
$ws.Cells.Item(50, 5).Value = "  +1.06%  "

$ws.Cells.Item(51, 5).Value = "  +3.14%  "
